$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.546146678731126
$ws.Range("C2").Value = 0.546146678731126
$ws.Range("D2").Value = 0.3706547201795419
$ws.Range("E2").Value = 0.6088141918348667
$ws.Range("F2").Value = 0.2784751697519821
$ws.Range("G2").Value = 15

$ws.Range("B3").Value = 0.3070856474631161
$ws.Range("C3").Value = 0.3070856474631161
$ws.Range("D3").Value = 0.1186079915109186
$ws.Range("E3").Value = 0.3443951095920478
$ws.Range("F3").Value = 0.1617903564833405
$ws.Range("G3").Value = 14

$ws.Range("B4").Value = 0.235791455158233
$ws.Range("C4").Value = 0.2409394054473722
$ws.Range("D4").Value = 0.08982400062291519
$ws.Range("E4").Value = 0.2997065241580757
$ws.Range("F4").Value = 0.1925580159035662
$ws.Range("G4").Value = 13

$ws.Range("B5").Value = 0.3438942213616346
$ws.Range("C5").Value = 0.3438942213616346
$ws.Range("D5").Value = 0.1608659895464369
$ws.Range("E5").Value = 0.4010810261610949
$ws.Range("F5").Value = 0.2155823084169401
$ws.Range("G5").Value = 12

$ws.Range("B6").Value = 0.3554346146561582
$ws.Range("C6").Value = 0.3554346146561582
$ws.Range("D6").Value = 0.166620911117176
$ws.Range("E6").Value = 0.4081922477426244
$ws.Range("F6").Value = 0.2105133259524082
$ws.Range("G6").Value = 11

$ws.Range("B7").Value = 0.3119807990786818
$ws.Range("C7").Value = 0.3119807990786818
$ws.Range("D7").Value = 0.1298725105494941
$ws.Range("E7").Value = 0.3603782881216544
$ws.Range("F7").Value = 0.1901475788133503
$ws.Range("G7").Value = 10

$ws.Range("B8").Value = 0.342631818604201
$ws.Range("C8").Value = 0.342631818604201
$ws.Range("D8").Value = 0.1511607671542049
$ws.Range("E8").Value = 0.3887939906354069
$ws.Range("F8").Value = 0.1948967150530138
$ws.Range("G8").Value = 9

$ws.Range("B9").Value = 0.3542225829761536
$ws.Range("C9").Value = 0.3542225829761536
$ws.Range("D9").Value = 0.162709113486236
$ws.Range("E9").Value = 0.4033721773824218
$ws.Range("F9").Value = 0.2062882177812338
$ws.Range("G9").Value = 8

$ws.Range("B10").Value = 0.3488257250224557
$ws.Range("C10").Value = 0.3488257250224557
$ws.Range("D10").Value = 0.1632182784925353
$ws.Range("E10").Value = 0.4040028199066626
$ws.Range("F10").Value = 0.2201409565050135
$ws.Range("G10").Value = 7

$ws.Range("B11").Value = 0.3139331037724223
$ws.Range("C11").Value = 0.3139331037724223
$ws.Range("D11").Value = 0.1244752993284096
$ws.Range("E11").Value = 0.352810571452174
$ws.Range("F11").Value = 0.1763677034523832
$ws.Range("G11").Value = 6
